# Applies the "cn181107" console upload edit to deleaveInneripDialog.xlsx
#
# Content changes (column C, the English translation column):
#   C4  : "Number of Available IPs in Subnet"   -> "Subnet IP Count"
#   C5  : "Available Quota/IP Quota"            -> "ENI IP Quota"
#   C12 : "Please add one secondary IP at least" -> "At lease one secondary IP is required."
#   C11 : "Please fill in the correct IP address" -> "Please enter a valid IP address."
#
# (C12 is written before C11 so the new shared-string table lands in the same
#  order as the authored workbook: …, "At lease one secondary IP is
#  required." then "Please enter a valid IP address.")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = "Subnet IP Count"
$ws.Range("C5").Value  = "ENI IP Quota"
$ws.Range("C12").Value = "At lease one secondary IP is required."
$ws.Range("C11").Value = "Please enter a valid IP address."

# Page setup: portrait, paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the selection / active cell to C19 (matches the saved view state)
$ws.Range("C19").Select()
